# Add the Mar-2021 "9626 HK Equity" liquidity override row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A5").Value = (Get-Date -Year 2021 -Month 3 -Day 31).Date
$ws.Range("B5").Value = "9626 HK Equity"
$ws.Range("C5").Value = "L0"

$ws.Range("B6").Select()
